$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.096.48'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.560.37'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.23%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.22'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.486'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.39%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.11'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.246'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.77%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0866'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.782.21'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.19%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.559.40'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.91%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.86%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '62.90'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.105.25'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '214.95'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0684'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.22'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.11%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.31'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -5.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.00'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.49'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.23%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -7.42%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.81%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.61%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.16'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.387.48'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.91'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.42%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.29'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.941'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -4.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0164'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.21%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.20%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.78'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.31'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.05%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.16'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.35%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.26'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.695.61'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.39'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₇0983'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.41%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.16%  '
